$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H116").Value = 5127147.5
$ws.Range("I116").Value = 6291727
$ws.Range("J116").Value = 2999.2
$ws.Range("K116").Value = 6291727
$ws.Range("L116").Value = 2999.2
$ws.Range("M116").Value = -6288285
$ws.Range("N116").Value = -9883.200000000001
$ws.Range("H133").Value = 29528.25
$ws.Range("J133").Value = 29528.25
$ws.Range("L133").Value = 29528.25
$ws.Range("N133").Value = -39648.25

$ws = $wb.Sheets("ARM")
$ws.Range("H32").Value = 16519.473
$ws.Range("I32").Value = 2622.677
$ws.Range("J32").Value = 145561.14
$ws.Range("K32").Value = 2622.677
$ws.Range("L32").Value = 145561.14
$ws.Range("M32").Value = -2335.677
$ws.Range("N32").Value = -146135.14
$ws.Range("H61").Value = 1817.525
$ws.Range("I61").Value = 1370.7812
$ws.Range("K61").Value = 1370.7812
$ws.Range("M61").Value = -1158.7812
$ws.Range("H132").Value = 2614.2927
$ws.Range("I132").Value = 2169.9644
$ws.Range("J132").Value = 3571.3076
$ws.Range("K132").Value = 6509.8932
$ws.Range("L132").Value = 10713.9228
$ws.Range("M132").Value = -3979.8932
$ws.Range("N132").Value = -15773.9228
$ws.Range("H136").Value = 1817.525
$ws.Range("I136").Value = 1370.7812
$ws.Range("K136").Value = 4112.3436
$ws.Range("M136").Value = -1562.3436

$ws = $wb.Sheets("BSM")
$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H134").Value = 23258542
$ws.Range("I134").Value = 35715620
$ws.Range("K134").Value = 107146860
$ws.Range("M134").Value = -107144325

$ws = $wb.Sheets("CRP")
$ws.Range("H31").Value = 1310.8975
$ws.Range("I31").Value = 968.7
$ws.Range("J31").Value = 2451.5557
$ws.Range("K31").Value = 968.7
$ws.Range("L31").Value = 2451.5557
$ws.Range("M31").Value = -673.7
$ws.Range("N31").Value = -3041.5557
$ws.Range("H34").Value = 1310.8975
$ws.Range("I34").Value = 968.7
$ws.Range("J34").Value = 2451.5557
$ws.Range("K34").Value = 968.7
$ws.Range("L34").Value = 2451.5557
$ws.Range("M34").Value = -766.7
$ws.Range("N34").Value = -2855.5557
$ws.Range("H58").Value = 2432.65
$ws.Range("I58").Value = 1053.1428
$ws.Range("J58").Value = 3175.4614
$ws.Range("K58").Value = 1053.1428
$ws.Range("L58").Value = 3175.4614
$ws.Range("M58").Value = -850.1428000000001
$ws.Range("N58").Value = -3581.4614
$ws.Range("H107").Value = 274.29166
$ws.Range("I107").Value = 190.125
$ws.Range("K107").Value = 190.125
$ws.Range("M107").Value = 1729.875
$ws.Range("H132").Value = 1867.3489
$ws.Range("I132").Value = 1158.909
$ws.Range("K132").Value = 3476.727
$ws.Range("M132").Value = -946.7270000000003
$ws.Range("H134").Value = 2390.532
$ws.Range("I134").Value = 1395.5588
$ws.Range("J134").Value = 4992.769
$ws.Range("K134").Value = 4186.6764
$ws.Range("L134").Value = 14978.307
$ws.Range("M134").Value = -1651.6764
$ws.Range("N134").Value = -20048.307
$ws.Range("H136").Value = 2432.65
$ws.Range("I136").Value = 1053.1428
$ws.Range("J136").Value = 3175.4614
$ws.Range("K136").Value = 3159.4284
$ws.Range("L136").Value = 9526.3842
$ws.Range("M136").Value = -609.4284000000002
$ws.Range("N136").Value = -14626.3842

$ws = $wb.Sheets("GSM")
$ws.Range("H96").Value = 24000
$ws.Range("J96").Value = 24000
$ws.Range("L96").Value = 24000
$ws.Range("N96").Value = -29492
$ws.Range("H102").Value = 2408.55
$ws.Range("I102").Value = 2096.7693
$ws.Range("J102").Value = 2987.5715
$ws.Range("K102").Value = 2096.7693
$ws.Range("L102").Value = 2987.5715
$ws.Range("M102").Value = -474.7692999999999
$ws.Range("N102").Value = -6231.5715
$ws.Range("H123").Value = 9301.333000000001
$ws.Range("J123").Value = 9301.333000000001
$ws.Range("L123").Value = 9301.333000000001
$ws.Range("N123").Value = -14201.333
$ws.Range("H132").Value = 2291.2659
$ws.Range("I132").Value = 2016.1936
$ws.Range("J132").Value = 3294.4707
$ws.Range("K132").Value = 6048.5808
$ws.Range("L132").Value = 9883.4121
$ws.Range("M132").Value = -3518.5808
$ws.Range("N132").Value = -14943.4121

$ws = $wb.Sheets("LTW")
$ws.Range("H22").Value = 9879.166999999999
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 9879.166999999999
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 9879.166999999999
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -10469.167
$ws.Range("H27").Value = 9879.166999999999
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 9879.166999999999
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 9879.166999999999
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -10093.167
$ws.Range("H46").Value = 1318.25
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1509.2
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1509.2
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1885.2
$ws.Range("H82").Value = 1212.1428
$ws.Range("I82").Value = 883.3333
$ws.Range("J82").Value = 1458.75
$ws.Range("K82").Value = 883.3333
$ws.Range("L82").Value = 1458.75
$ws.Range("M82").Value = -522.3333
$ws.Range("N82").Value = -2180.75
$ws.Range("H85").Value = 1212.1428
$ws.Range("I85").Value = 883.3333
$ws.Range("J85").Value = 1458.75
$ws.Range("K85").Value = 883.3333
$ws.Range("L85").Value = 1458.75
$ws.Range("M85").Value = 364.6667
$ws.Range("N85").Value = -3954.75
$ws.Range("H93").Value = 1111.95
$ws.Range("I93").Value = 849.6429000000001
$ws.Range("J93").Value = 1724
$ws.Range("K93").Value = 849.6429000000001
$ws.Range("L93").Value = 1724
$ws.Range("M93").Value = 398.3570999999999
$ws.Range("N93").Value = -4220
$ws.Range("H94").Value = 19000
$ws.Range("J94").Value = 19000
$ws.Range("L94").Value = 19000
$ws.Range("N94").Value = -20352
$ws.Range("H100").Value = 2978919.2
$ws.Range("I100").Value = 10418701
$ws.Range("J100").Value = 3006.6667
$ws.Range("K100").Value = 10418701
$ws.Range("L100").Value = 3006.6667
$ws.Range("M100").Value = -10418160
$ws.Range("N100").Value = -4088.6667
$ws.Range("H122").Value = 3195.8518
$ws.Range("I122").Value = 2424
$ws.Range("J122").Value = 3813.3333
$ws.Range("K122").Value = 7272
$ws.Range("L122").Value = 11439.9999
$ws.Range("M122").Value = -4822
$ws.Range("N122").Value = -16339.9999

$ws = $wb.Sheets("WVR")
$ws.Range("H95").Value = 500344
$ws.Range("J95").Value = 500344
$ws.Range("L95").Value = 500344
$ws.Range("N95").Value = -505836
$ws.Range("H107").Value = 2315743.8
$ws.Range("I107").Value = 3087294.5
$ws.Range("K107").Value = 9261883.5
$ws.Range("M107").Value = -9259963.5
$ws.Range("H132").Value = 10640563
$ws.Range("I132").Value = 12822213
$ws.Range("J132").Value = 5018
$ws.Range("K132").Value = 38466639
$ws.Range("L132").Value = 15054
$ws.Range("M132").Value = -38464109
$ws.Range("N132").Value = -20114
